$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row above the current row 15 ("PathPDriveFolder"), shifting
# that row (and everything below it) down by one.
$ws.Rows(15).Insert()

# Populate the newly inserted row with the new "different P drive folder"
# setting.
$ws.Range("A15").Value = "PathDifferentPDriveFolder"
$ws.Range("B15").Value = "\\10.250.52.158\Depts\TaxReturnOutSourcing\Preparer\UIPathPublish\IR Bot Temp Files\InputFiles"
$ws.Range("C15").Value = "\\somproddfs1.prod.sovos.org\depts\TaxSolver Files"

# Match the row height / formatting of the neighboring config rows.
$ws.Rows(15).RowHeight = 14.25
$ws.Range("A15:D15").Interior.Color = 5296274

# Make "Settings" the active sheet/tab and put the selection on the new row,
# matching where the author was working when they saved.
$ws.Activate()
$ws.Range("A15").Select()
